# Append patient rows 2-10 to the active sheet, matching the source data.
# Values that look numeric (ID, Age, TelegramID) are forced to Text via a
# leading apostrophe, preserving the original workbook's convention of
# storing every cell as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("3",  "Әсем Ерболқызы",    "23", "Невролог", "6863633722"),
    @("4",  "Жанар Амангелді",   "17", "Хирург",   "6863633722"),
    @("5",  "Гүлім Айтжан",      "21", "Хирург",   "6863633722"),
    @("6",  "Аяна Бақыт",        "15", "Терапевт", "6863633722"),
    @("7",  "Сабина Жеңіс",      "23", "Невролог", "6863633722"),
    @("8",  "Айгерім Нұрмұхан",  "22", "Хирург",   "6863633722"),
    @("1",  "Аружан Кәрім",      "23", "Хирург",   "6863633722"),
    @("9",  "Алина Төлеген",     "56", "Терапевт", "6863633722"),
    @("10", "Дильнара Сағындық", "56", "Невролог", "6863633722")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = "'" + $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = "'" + $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = "'" + $row[4]
    $rowIndex++
}
